$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = [double]"0.9922803468680383"
$ws.Range("E2").Value = [double]"0.9922803468680383"

# Row 3
$ws.Range("D3").Value = [double]"1.376189089261484E-26"
$ws.Range("E3").Value = [double]"1.376189089261484E-26"

# Row 4
$ws.Range("D4").Value = [double]"5.543616232827255E-19"
$ws.Range("E4").Value = [double]"5.543616232827255E-19"

# Row 5
$ws.Range("D5").Value = [double]"0.9999999999994191"
$ws.Range("E5").Value = [double]"0.9999999999994191"

# Row 6
$ws.Range("C6").Value = $false
$ws.Range("D6").Value = [double]"0.9992481731165173"
$ws.Range("E6").Value = [double]"0.9992481731165173"

# Row 7
$ws.Range("D7").Value = [double]"8.055981314646563E-14"
$ws.Range("E7").Value = [double]"0.9999999999999194"

# Row 8
$ws.Range("D8").Value = [double]"0.9999999999932081"
$ws.Range("E8").Value = [double]"6.791900375446858E-12"

# Row 9
$ws.Range("D9").Value = [double]"0.9992138498385643"
$ws.Range("E9").Value = [double]"0.0007861501614356836"

# Row 10
$ws.Range("D10").Value = [double]"5.392651757992452E-24"

# Row 11
$ws.Range("D11").Value = [double]"1.034458475777952E-16"
$ws.Range("E11").Value = [double]"0.9999999999999999"
$ws.Range("F11").Value = [double]"16.07662963867188"
$ws.Range("G11").Value = [double]"0.4"
